# Update the NATMI Vcam1-Itga4 LR-pair results sheet with the
# re-run values (3 replicates) and the full 4x4 cluster combinations
# (ECs / FAPs / M2 / sCs), per Dr Hou's advice.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'ECs'
$ws.Range("B2").Value = 'Vcam1'
$ws.Range("C2").Value = 'Itga4'
$ws.Range("D2").Value = 'ECs'
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.312497333333334
$ws.Range("H2").Value = 21.937492
$ws.Range("I2").Value = 0.05970572560549242
$ws.Range("J2").Value = 0.05970572560549242
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 23.65990166666667
$ws.Range("N2").Value = 70.979705
$ws.Range("O2").Value = 0.2997993941754699
$ws.Range("P2").Value = 0.29979939417547
$ws.Range("Q2").Value = 173.0129678444289
$ws.Range("R2").Value = 1557.11671059986
$ws.Range("S2").Value = 0.01789974036533347
$ws.Range("T2").Value = 0.01789974036533347

# Row 3
$ws.Range("A3").Value = 'ECs'
$ws.Range("B3").Value = 'Vcam1'
$ws.Range("C3").Value = 'Itga4'
$ws.Range("D3").Value = 'FAPs'
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.312497333333334
$ws.Range("H3").Value = 21.937492
$ws.Range("I3").Value = 0.05970572560549242
$ws.Range("J3").Value = 0.05970572560549242
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.07690566666666666
$ws.Range("N3").Value = 0.230717
$ws.Range("O3").Value = 0.0009744872400636476
$ws.Range("P3").Value = 0.0009744872400636479
$ws.Range("Q3").Value = 0.5623724824182222
$ws.Range("R3").Value = 5.061352341764001
$ws.Range("S3").Value = 0.00005818246776129376
$ws.Range("T3").Value = 0.00005818246776129378

# Row 4
$ws.Range("A4").Value = 'ECs'
$ws.Range("B4").Value = 'Vcam1'
$ws.Range("C4").Value = 'Itga4'
$ws.Range("D4").Value = 'M2'
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.312497333333334
$ws.Range("H4").Value = 21.937492
$ws.Range("I4").Value = 0.05970572560549242
$ws.Range("J4").Value = 0.05970572560549242
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 53.21452433333334
$ws.Range("N4").Value = 159.643573
$ws.Range("O4").Value = 0.6742919890890982
$ws.Range("P4").Value = 0.6742919890890983
$ws.Range("Q4").Value = 389.1310672821018
$ws.Range("R4").Value = 3502.179605538916
$ws.Range("S4").Value = 0.04025909247853538
$ws.Range("T4").Value = 0.04025909247853539

# Row 5
$ws.Range("A5").Value = 'ECs'
$ws.Range("B5").Value = 'Vcam1'
$ws.Range("C5").Value = 'Itga4'
$ws.Range("D5").Value = 'sCs'
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.312497333333334
$ws.Range("H5").Value = 21.937492
$ws.Range("I5").Value = 0.05970572560549242
$ws.Range("J5").Value = 0.05970572560549242
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.967779333333333
$ws.Range("N5").Value = 5.903338
$ws.Range("O5").Value = 0.02493412949536815
$ws.Range("P5").Value = 0.02493412949536816
$ws.Range("Q5").Value = 14.38938112758844
$ws.Range("R5").Value = 129.504430148296
$ws.Range("S5").Value = 0.001488710293862266
$ws.Range("T5").Value = 0.001488710293862266

# Row 6
$ws.Range("A6").Value = 'FAPs'
$ws.Range("B6").Value = 'Vcam1'
$ws.Range("C6").Value = 'Itga4'
$ws.Range("D6").Value = 'ECs'
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 29.68221266666667
$ws.Range("H6").Value = 89.046638
$ws.Range("I6").Value = 0.242351957758873
$ws.Range("J6").Value = 0.242351957758873
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 23.65990166666667
$ws.Range("N6").Value = 70.979705
$ws.Range("O6").Value = 0.2997993941754699
$ws.Range("P6").Value = 0.29979939417547
$ws.Range("Q6").Value = 702.2782329424211
$ws.Range("R6").Value = 6320.50409648179
$ws.Range("S6").Value = 0.0726569701133492
$ws.Range("T6").Value = 0.0726569701133492

# Row 7
$ws.Range("A7").Value = 'FAPs'
$ws.Range("B7").Value = 'Vcam1'
$ws.Range("C7").Value = 'Itga4'
$ws.Range("D7").Value = 'FAPs'
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 29.68221266666667
$ws.Range("H7").Value = 89.046638
$ws.Range("I7").Value = 0.242351957758873
$ws.Range("J7").Value = 0.242351957758873
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.07690566666666666
$ws.Range("N7").Value = 0.230717
$ws.Range("O7").Value = 0.0009744872400636476
$ws.Range("P7").Value = 0.0009744872400636479
$ws.Range("Q7").Value = 2.282730353271778
$ws.Range("R7").Value = 20.544573179446
$ws.Range("S7").Value = 0.0002361688904404659
$ws.Range("T7").Value = 0.0002361688904404659

# Row 8
$ws.Range("A8").Value = 'FAPs'
$ws.Range("B8").Value = 'Vcam1'
$ws.Range("C8").Value = 'Itga4'
$ws.Range("D8").Value = 'M2'
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 29.68221266666667
$ws.Range("H8").Value = 89.046638
$ws.Range("I8").Value = 0.242351957758873
$ws.Range("J8").Value = 0.242351957758873
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 53.21452433333334
$ws.Range("N8").Value = 159.643573
$ws.Range("O8").Value = 0.6742919890890982
$ws.Range("P8").Value = 0.6742919890890983
$ws.Range("Q8").Value = 1579.524828217508
$ws.Range("R8").Value = 14215.72345395757
$ws.Range("S8").Value = 0.1634159836568676
$ws.Range("T8").Value = 0.1634159836568676

# Row 9
$ws.Range("A9").Value = 'FAPs'
$ws.Range("B9").Value = 'Vcam1'
$ws.Range("C9").Value = 'Itga4'
$ws.Range("D9").Value = 'sCs'
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 29.68221266666667
$ws.Range("H9").Value = 89.046638
$ws.Range("I9").Value = 0.242351957758873
$ws.Range("J9").Value = 0.242351957758873
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.967779333333333
$ws.Range("N9").Value = 5.903338
$ws.Range("O9").Value = 0.02493412949536815
$ws.Range("P9").Value = 0.02493412949536816
$ws.Range("Q9").Value = 58.40804465307156
$ws.Range("R9").Value = 525.672401877644
$ws.Range("S9").Value = 0.006042835098215731
$ws.Range("T9").Value = 0.006042835098215732

# Row 10
$ws.Range("A10").Value = 'M2'
$ws.Range("B10").Value = 'Vcam1'
$ws.Range("C10").Value = 'Itga4'
$ws.Range("D10").Value = 'ECs'
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 11.06470466666667
$ws.Range("H10").Value = 33.194114
$ws.Range("I10").Value = 0.09034208022509747
$ws.Range("J10").Value = 0.09034208022509749
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 23.65990166666667
$ws.Range("N10").Value = 70.979705
$ws.Range("O10").Value = 0.2997993941754699
$ws.Range("P10").Value = 0.29979939417547
$ws.Range("Q10").Value = 261.7898243840411
$ws.Range("R10").Value = 2356.10841945637
$ws.Range("S10").Value = 0.02708450092003592
$ws.Range("T10").Value = 0.02708450092003593

# Row 11
$ws.Range("A11").Value = 'M2'
$ws.Range("B11").Value = 'Vcam1'
$ws.Range("C11").Value = 'Itga4'
$ws.Range("D11").Value = 'FAPs'
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 11.06470466666667
$ws.Range("H11").Value = 33.194114
$ws.Range("I11").Value = 0.09034208022509747
$ws.Range("J11").Value = 0.09034208022509749
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.07690566666666666
$ws.Range("N11").Value = 0.230717
$ws.Range("O11").Value = 0.0009744872400636476
$ws.Range("P11").Value = 0.0009744872400636479
$ws.Range("Q11").Value = 0.8509384888597776
$ws.Range("R11").Value = 7.658446399738
$ws.Range("S11").Value = 0.00008803720442016387
$ws.Range("T11").Value = 0.00008803720442016391

# Row 12
$ws.Range("A12").Value = 'M2'
$ws.Range("B12").Value = 'Vcam1'
$ws.Range("C12").Value = 'Itga4'
$ws.Range("D12").Value = 'M2'
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 11.06470466666667
$ws.Range("H12").Value = 33.194114
$ws.Range("I12").Value = 0.09034208022509747
$ws.Range("J12").Value = 0.09034208022509749
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 53.21452433333334
$ws.Range("N12").Value = 159.643573
$ws.Range("O12").Value = 0.6742919890890982
$ws.Range("P12").Value = 0.6742919890890983
$ws.Range("Q12").Value = 588.8029957254802
$ws.Range("R12").Value = 5299.226961529322
$ws.Range("S12").Value = 0.06091694097342785
$ws.Range("T12").Value = 0.06091694097342788

# Row 13
$ws.Range("A13").Value = 'M2'
$ws.Range("B13").Value = 'Vcam1'
$ws.Range("C13").Value = 'Itga4'
$ws.Range("D13").Value = 'sCs'
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 11.06470466666667
$ws.Range("H13").Value = 33.194114
$ws.Range("I13").Value = 0.09034208022509747
$ws.Range("J13").Value = 0.09034208022509749
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.967779333333333
$ws.Range("N13").Value = 5.903338
$ws.Range("O13").Value = 0.02493412949536815
$ws.Range("P13").Value = 0.02493412949536816
$ws.Range("Q13").Value = 21.77289717250355
$ws.Range("R13").Value = 195.956074552532
$ws.Range("S13").Value = 0.002252601127213518
$ws.Range("T13").Value = 0.00225260112721352

# Row 14
$ws.Range("A14").Value = 'sCs'
$ws.Range("B14").Value = 'Vcam1'
$ws.Range("C14").Value = 'Itga4'
$ws.Range("D14").Value = 'ECs'
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 74.41623166666666
$ws.Range("H14").Value = 223.248695
$ws.Range("I14").Value = 0.6076002364105371
$ws.Range("J14").Value = 0.6076002364105371
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 23.65990166666667
$ws.Range("N14").Value = 70.979705
$ws.Range("O14").Value = 0.2997993941754699
$ws.Range("P14").Value = 0.29979939417547
$ws.Range("Q14").Value = 1760.680723637219
$ws.Range("R14").Value = 15846.12651273497
$ws.Range("S14").Value = 0.1821581827767513
$ws.Range("T14").Value = 0.1821581827767513

# Row 15
$ws.Range("A15").Value = 'sCs'
$ws.Range("B15").Value = 'Vcam1'
$ws.Range("C15").Value = 'Itga4'
$ws.Range("D15").Value = 'FAPs'
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 74.41623166666666
$ws.Range("H15").Value = 223.248695
$ws.Range("I15").Value = 0.6076002364105371
$ws.Range("J15").Value = 0.6076002364105371
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.07690566666666666
$ws.Range("N15").Value = 0.230717
$ws.Range("O15").Value = 0.0009744872400636476
$ws.Range("P15").Value = 0.0009744872400636479
$ws.Range("Q15").Value = 5.723029907146111
$ws.Range("R15").Value = 51.507269164315
$ws.Range("S15").Value = 0.000592098677441724
$ws.Range("T15").Value = 0.0005920986774417242

# Row 16
$ws.Range("A16").Value = 'sCs'
$ws.Range("B16").Value = 'Vcam1'
$ws.Range("C16").Value = 'Itga4'
$ws.Range("D16").Value = 'M2'
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 74.41623166666666
$ws.Range("H16").Value = 223.248695
$ws.Range("I16").Value = 0.6076002364105371
$ws.Range("J16").Value = 0.6076002364105371
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 53.21452433333334
$ws.Range("N16").Value = 159.643573
$ws.Range("O16").Value = 0.6742919890890982
$ws.Range("P16").Value = 0.6742919890890983
$ws.Range("Q16").Value = 3960.024370820804
$ws.Range("R16").Value = 35640.21933738724
$ws.Range("S16").Value = 0.4096999719802673
$ws.Range("T16").Value = 0.4096999719802674

# Row 17
$ws.Range("A17").Value = 'sCs'
$ws.Range("B17").Value = 'Vcam1'
$ws.Range("C17").Value = 'Itga4'
$ws.Range("D17").Value = 'sCs'
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 74.41623166666666
$ws.Range("H17").Value = 223.248695
$ws.Range("I17").Value = 0.6076002364105371
$ws.Range("J17").Value = 0.6076002364105371
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.967779333333333
$ws.Range("N17").Value = 5.903338
$ws.Range("O17").Value = 0.02493412949536815
$ws.Range("P17").Value = 0.02493412949536816
$ws.Range("Q17").Value = 146.4347227382122
$ws.Range("R17").Value = 1317.91250464391
$ws.Range("S17").Value = 0.01514998297607663
$ws.Range("T17").Value = 0.01514998297607664

Write-Output "Updated rows 2-17 (A:T) with re-run NATMI results"
